# Apply updated crypto price/volume figures (scheduled GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.183.48'
$ws.Range('E2').Value = '  -1.21%  '
$ws.Range('D3').Value = '1.894.82'
$ws.Range('E3').Value = '  -0.90%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '''245.60'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('D6').Value = '''0.685'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +8.13%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').Value = '''40.42'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.59%  '
$ws.Range('D9').Value = '''0.345'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.83%  '
$ws.Range('D10').Value = '''53.06'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +10.84%  '
$ws.Range('D11').Value = '''0.0717'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Value = '''0.0988'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.82%  '
$ws.Range('D13').Value = '2.169.15'
$ws.Range('E13').Value = '  -0.87%  '
$ws.Range('D14').Value = '''12.50'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.12%  '
$ws.Range('D15').Value = '''0.701'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.05%  '
$ws.Range('D16').Value = '1.890.85'
$ws.Range('E16').Value = '  -1.49%  '
$ws.Range('E17').Value = '  -1.66%  '
$ws.Range('D18').Value = '35.202.74'
$ws.Range('E18').Value = '  -1.10%  '
$ws.Range('D19').Value = '''71.86'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.17%  '
$ws.Range('E20').Value = '  +0.39%  '
$ws.Range('D21').Value = '''239.88'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.65%  '
$ws.Range('D22').Value = '''12.53'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.26%  '
$ws.Range('D23').Value = '''4.76'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.67%  '
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('E25').Value = '  +1.13%  '
$ws.Range('D26').Value = '''2.33'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +10.31%  '
$ws.Range('D27').Value = '''167.82'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Value = '''8.50'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.22%  '
$ws.Range('D29').Value = '''0.130'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.42%  '
$ws.Range('D30').Value = '''18.17'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.87%  '
$ws.Range('E32').Value = '  +0.75%  '
$ws.Range('D33').Value = '''0.0562'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.69%  '
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('D35').Value = '''1.85'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +7.43%  '
$ws.Range('E36').Value = '  -2.26%  '
$ws.Range('D37').Value = '''0.902'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.29%  '
$ws.Range('D38').Value = '''1.48'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +11.61%  '
$ws.Range('E39').Value = '  -1.82%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').Value = '''1.09'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.43%  '
$ws.Range('E41').Value = '  +0.94%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = '''0.0642'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +8.00%  '
$ws.Range('D43').Value = '''15.99'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.98%  '
$ws.Range('D44').Value = '''89.25'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.78%  '
$ws.Range('D45').Value = '1.338.97'
$ws.Range('E45').Value = '  -1.61%  '
$ws.Range('D46').Value = '''2.40'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.81%  '
$ws.Range('E47').Value = '  -0.31%  '
$ws.Range('E48').Value = '  -0.09%  '
$ws.Range('D49').Value = '''45.05'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -8.50%  '
$ws.Range('D50').Value = '''12.17'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.04%  '
$ws.Range('E51').Value = '  -3.60%  '
